# Insert two new columns ("Limite (dBµV/m)" and "Margin (dB)") right before
# the existing "Polarization" column (H). Excel shifts H:N right to J:P.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H:I").Insert()

# The insert leaves the sheet spanning A:P (old M/N now duplicated out at O/P).
# Every cell from H to N is rewritten explicitly below to the final target
# values, so just drop the now-unused overflow columns O:P afterwards.

# ---- Header row (row 1) ----
$ws.Range("H1").Value = "Limite (dBµV/m)"
$ws.Range("I1").Value = "Margin (dB)"
$ws.Range("J1").Value = "Polarization"
$ws.Range("K1").Value = "Correction (dB)"
$ws.Range("L1").Value = "Overtaking (dB)"
$ws.Range("M1").Value = "Conformity"
$ws.Range("N1").Value = "Configuration"

# ---- Data rows ----
# Row 2 (Peak, SR1)
$ws.Range("F2").Value = 33.25
$ws.Range("H2").Value = 64
$ws.Range("I2").Value = 30.75
$ws.Range("J2").Value = "Vertical"
$ws.Range("K2").Value = 9.92
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = "OK"
$ws.Range("N2").Value = "ER_In front of harness RBW 9kHz"

# Row 3 (Peak, SR2)
$ws.Range("F3").Value = 22.15
$ws.Range("H3").Value = 58
$ws.Range("I3").Value = 35.85
$ws.Range("J3").Value = "Vertical"
$ws.Range("K3").Value = 10.29
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = "OK"
$ws.Range("N3").Value = "ER_In front of harness RBW 9kHz"

# Row 4 (Section/Detector type relabeled Peak -> Q-Peak, SR1)
$ws.Range("C4").Value = "Q-Peak"
$ws.Range("F4").Value = 31.39
$ws.Range("G4").Value = "Q-Peak"
$ws.Range("H4").Value = 51
$ws.Range("I4").Value = 19.61
$ws.Range("J4").Value = "Vertical"
$ws.Range("K4").Value = 9.92
$ws.Range("L4").Value = "-"
$ws.Range("M4").Value = "OK"
$ws.Range("N4").Value = "ER_In front of harness RBW 9kHz"

# Row 5 (Peak -> Q-Peak, SR2)
$ws.Range("C5").Value = "Q-Peak"
$ws.Range("F5").Value = 16.05
$ws.Range("G5").Value = "Q-Peak"
$ws.Range("H5").Value = 45
$ws.Range("I5").Value = 28.95
$ws.Range("J5").Value = "Vertical"
$ws.Range("K5").Value = 10.29
$ws.Range("L5").Value = "-"
$ws.Range("M5").Value = "OK"
$ws.Range("N5").Value = "ER_In front of harness RBW 9kHz"

# Row 6 (Peak, SR1)
$ws.Range("F6").Value = 33.25
$ws.Range("H6").Value = 51
$ws.Range("I6").Value = 17.75
$ws.Range("J6").Value = "Vertical"
$ws.Range("K6").Value = 9.92
$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = "OK"
$ws.Range("N6").Value = "ER_In front of harness RBW 9kHz"

# Row 7 (Peak, SR2)
$ws.Range("F7").Value = 22.15
$ws.Range("H7").Value = 45
$ws.Range("I7").Value = 22.85
$ws.Range("J7").Value = "Vertical"
$ws.Range("K7").Value = 10.29
$ws.Range("L7").Value = "-"
$ws.Range("M7").Value = "OK"
$ws.Range("N7").Value = "ER_In front of harness RBW 9kHz"

# Row 8 (CISPR.AVG, SR1)
$ws.Range("F8").Value = 18.35
$ws.Range("H8").Value = 44
$ws.Range("I8").Value = 25.65
$ws.Range("J8").Value = "Vertical"
$ws.Range("K8").Value = 9.92
$ws.Range("L8").Value = "-"
$ws.Range("M8").Value = "OK"
$ws.Range("N8").Value = "ER_In front of harness RBW 9kHz"

# Row 9 (CISPR.AVG, SR2)
$ws.Range("F9").Value = 4.57
$ws.Range("H9").Value = 28
$ws.Range("I9").Value = 23.43
$ws.Range("J9").Value = "Vertical"
$ws.Range("K9").Value = 10.29
$ws.Range("L9").Value = "-"
$ws.Range("M9").Value = "OK"
$ws.Range("N9").Value = "ER_In front of harness RBW 9kHz"

# Drop the overflow columns that the 2-column insert pushed past N.
$ws.Columns("O:P").Delete()
